# Edit script: adds D1-D4 "Respuestas" answer-key table (columns F:G, rows 2-46)
# mirroring the existing A/Respuesta/Respuestas table (B:D), plus Table2 ListObject.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Column widths for the new F:G block ---
$ws.Range("F1").ColumnWidth = 10.14
$ws.Range("G1").ColumnWidth = 10.14

# --- Header row ---
$ws.Range("F2").Value = "D"
$ws.Range("G2").Value = "Respuestas"

# --- Data rows 3-46 ---
$ws.Range("F3").Value = 0.1
$ws.Range("F4").Value = 0.2
$ws.Range("G4").Value = "C"
$ws.Range("F5").Value = 0.3
$ws.Range("F6").Value = 0.4
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = "D"
$ws.Range("F8").Value = 0.6
$ws.Range("G8").Value = "B"
$ws.Range("F9").Value = 0.7
$ws.Range("G9").Value = "C"
$ws.Range("F10").Value = 0.8
$ws.Range("G10").Value = "D"
$ws.Range("F11").Value = 0.9
$ws.Range("G11").Value = "C"
$ws.Range("F12").Value = 0.1
$ws.Range("G12").Value = "A"
$ws.Range("F13").Value = 0.11
$ws.Range("F14").Value = 2
$ws.Range("F15").Value = 0.1
$ws.Range("G15").Value = "B"
$ws.Range("F16").Value = 0.2
$ws.Range("G16").Value = "C"
$ws.Range("F17").Value = 0.3
$ws.Range("G17").Value = "B"
$ws.Range("F18").Value = 0.4
$ws.Range("G18").Value = "C"
$ws.Range("F19").Value = 0.5
$ws.Range("G19").Value = "B"
$ws.Range("F20").Value = 0.6
$ws.Range("G20").Value = "A"
$ws.Range("F21").Value = 0.7
$ws.Range("G21").Value = "D"
$ws.Range("F22").Value = 0.8
$ws.Range("G22").Value = "B"
$ws.Range("F23").Value = 0.9
$ws.Range("G23").Value = "B"
$ws.Range("F24").Value = 0.1
$ws.Range("G24").Value = "A"
$ws.Range("F25").Value = 3
$ws.Range("F26").Value = 0.1
$ws.Range("G26").Value = "B"
$ws.Range("F27").Value = 0.2
$ws.Range("G27").Value = "D"
$ws.Range("F28").Value = 0.3
$ws.Range("G28").Value = "B"
$ws.Range("F29").Value = 0.4
$ws.Range("G29").Value = "B"
$ws.Range("F30").Value = 0.5
$ws.Range("G30").Value = "A"
$ws.Range("F31").Value = 0.6
$ws.Range("G31").Value = "D"
$ws.Range("F32").Value = 0.7
$ws.Range("G32").Value = "A"
$ws.Range("F33").Value = 0.8
$ws.Range("G33").Value = "B"
$ws.Range("F34").Value = 0.9
$ws.Range("G34").Value = "C"
$ws.Range("F35").Value = 0.1
$ws.Range("G35").Value = "D"
$ws.Range("F36").Value = 4
$ws.Range("F37").Value = 0.1
$ws.Range("G37").Value = "D"
$ws.Range("F38").Value = 0.2
$ws.Range("G38").Value = "C"
$ws.Range("F39").Value = 0.3
$ws.Range("G39").Value = "D"
$ws.Range("F40").Value = 0.4
$ws.Range("G40").Value = "A"
$ws.Range("F41").Value = 0.5
$ws.Range("G41").Value = "D"
$ws.Range("F42").Value = 0.6
$ws.Range("G42").Value = "C"
$ws.Range("F43").Value = 0.7
$ws.Range("G43").Value = "C"
$ws.Range("F44").Value = 0.8
$ws.Range("G44").Value = "B"
$ws.Range("F45").Value = 0.9
$ws.Range("G45").Value = "A"
$ws.Range("F46").Value = 0.1
$ws.Range("G46").Value = "A"

# --- Apply "0.00" number format to the 4 "D1.10/D2.10/D3.10/D4.10" style cells ---
# (these round decade markers get numFmtId 2 in the target, same as the diff)
$ws.Range("F12").NumberFormat = "0.00"
$ws.Range("F24").NumberFormat = "0.00"
$ws.Range("F35").NumberFormat = "0.00"
$ws.Range("F46").NumberFormat = "0.00"

# --- Create the Table2 ListObject over F2:G46 ---
$tbl2 = $ws.ListObjects.Add(1, $ws.Range("F2:G46"), 0, 1)
$tbl2.Name = "Table2"
$tbl2.TableStyle = "TableStyleLight1"

# --- Update sheet view: scroll position + selection ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 24
[void]$ws.Range("G47").Select()

Write-Output "D1-D4 answer key table added"
